$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: COX4 -> MT1 (sequence unchanged)
$ws.Range("A6").Value = 'MT1'

# Insert a new row at position 9 (shifts old rows 9-11 down to 10-12),
# making room for the new ER1 entry.
$ws.Rows.Item(9).Insert()

# New row 9: ER1
$ws.Range("A9").Value = 'ER1'
$ws.Range("B9").Value = 'ATATTAGAGCAACCTCTGAAAT TTGTGCTTACTGCGGCCGTCGTG CTCTTGACGACGTCGGTTCTTTG TTGTGTAGTATTTACA'
$ws.Range("B9").Style = "Normal"

# Former row 9 (CYB5) is now row 10 -> rename tag to ER2, sequence kept
$ws.Range("A10").Value = 'ER2'

# Former row 10 (SNC1) is now row 11 -> rename tag to PM1, sequence kept
$ws.Range("A11").Value = 'PM1'

# Former row 11 (NES1) is now row 12, unchanged

# New row 13: VC1
$ws.Range("A13").Value = 'VC1'
$ws.Range("B13").Value = 'AATATAAAAGAAATAATGTGGT GGCAGAAGGTCAAAAATATTAC GTTATTAACTTTCACTATTATAC TATTTGTAAGTGCTGCTTTCATG TTTTTCTATCTGTGG'
$ws.Range("B13").Style = "Normal"

$ws.Range("A2:A13").Select()
